# Append the latest EUR->ARS quote as a new row at the bottom of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 65

# Column A holds a date-shaped string ("2025-10-08"); without a text hint
# Excel would silently reinterpret it as a serial date. A leading apostrophe
# (the same trick a person typing it into Excel would use) keeps it as a
# literal string, matching the existing rows above it.
$ws.Cells.Item($newRow, 1).Value = "'2025-10-08"
$ws.Cells.Item($newRow, 2).Value = "15:22:59"
$ws.Cells.Item($newRow, 3).Value = "1.00 EUR = 1,770.2348"
